# Auto-generated script applying scheduled-runner price/profit refresh
# to the Leve profit tracking tables across all job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1293.6897
$ws.Range("I15").Value = 1293.6897
$ws.Range("K15").Value = 3881.0691
$ws.Range("M15").Value = -3712.0691
# Row 17
$ws.Range("H17").Value = 1372.5
$ws.Range("J17").Value = 1372.5
$ws.Range("L17").Value = 4117.5
$ws.Range("N17").Value = -4453.5
# Row 19
$ws.Range("H19").Value = 1149.4286
$ws.Range("I19").Value = 601.25
$ws.Range("J19").Value = 1880.3334
$ws.Range("K19").Value = 601.25
$ws.Range("L19").Value = 1880.3334
$ws.Range("M19").Value = -426.25
$ws.Range("N19").Value = -2230.3334
# Row 33
$ws.Range("H33").Value = 204.375
$ws.Range("I33").Value = 204.375
$ws.Range("K33").Value = 204.375
$ws.Range("M33").Value = 24.625
# Row 88
$ws.Range("H88").Value = 9849.166999999999
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
# Row 91
$ws.Range("H91").Value = 9849.166999999999
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
# Row 107
$ws.Range("H107").Value = 292
$ws.Range("I107").Value = 292
$ws.Range("K107").Value = 292
$ws.Range("M107").Value = 1628
# Row 138
$ws.Range("H138").Value = 4489.4546
$ws.Range("I138").Value = 2171
$ws.Range("J138").Value = 5814.2856
$ws.Range("K138").Value = 6513
$ws.Range("L138").Value = 17442.8568
$ws.Range("M138").Value = -1373
$ws.Range("N138").Value = -27722.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3503701.5
$ws.Range("I32").Value = 5838450
$ws.Range("J32").Value = 702003.3
$ws.Range("K32").Value = 5838450
$ws.Range("L32").Value = 702003.3
$ws.Range("M32").Value = -5838163
$ws.Range("N32").Value = -702577.3
# Row 132
$ws.Range("H132").Value = 2804.5
$ws.Range("I132").Value = 2916.5715
$ws.Range("K132").Value = 8749.7145
$ws.Range("M132").Value = -6219.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 823.25
$ws.Range("J64").Value = 749.2
$ws.Range("L64").Value = 749.2
$ws.Range("N64").Value = -1199.2
# Row 67
$ws.Range("H67").Value = 823.25
$ws.Range("J67").Value = 749.2
$ws.Range("L67").Value = 749.2
$ws.Range("N67").Value = -2309.2
# Row 76
$ws.Range("H76").Value = 200313.5
$ws.Range("J76").Value = 200313.5
$ws.Range("L76").Value = 200313.5
$ws.Range("N76").Value = -200943.5
# Row 79
$ws.Range("H79").Value = 200313.5
$ws.Range("J79").Value = 200313.5
$ws.Range("L79").Value = 200313.5
$ws.Range("N79").Value = -202497.5
# Row 88
$ws.Range("H88").Value = 20633.334
$ws.Range("J88").Value = 22760
$ws.Range("L88").Value = 22760
$ws.Range("N88").Value = -23572
# Row 91
$ws.Range("H91").Value = 20633.334
$ws.Range("J91").Value = 22760
$ws.Range("L91").Value = 22760
$ws.Range("N91").Value = -25568
# Row 94
$ws.Range("H94").Value = 976.8
$ws.Range("I94").Value = 1295
$ws.Range("J94").Value = 499.5
$ws.Range("K94").Value = 1295
$ws.Range("L94").Value = 499.5
$ws.Range("M94").Value = -844
$ws.Range("N94").Value = -1401.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 3642.7144
$ws.Range("J4").Value = 5249.5
$ws.Range("L4").Value = 5249.5
$ws.Range("N4").Value = -5473.5
# Row 7
$ws.Range("H7").Value = 75.84614999999999
$ws.Range("I7").Value = 46.5
$ws.Range("K7").Value = 46.5
$ws.Range("M7").Value = 66.5
# Row 31
$ws.Range("H31").Value = 1066.7
$ws.Range("I31").Value = 1013.4
$ws.Range("J31").Value = 1120
$ws.Range("K31").Value = 1013.4
$ws.Range("L31").Value = 1120
$ws.Range("M31").Value = -718.4
$ws.Range("N31").Value = -1710
# Row 34
$ws.Range("H34").Value = 1066.7
$ws.Range("I34").Value = 1013.4
$ws.Range("J34").Value = 1120
$ws.Range("K34").Value = 1013.4
$ws.Range("L34").Value = 1120
$ws.Range("M34").Value = -811.4
$ws.Range("N34").Value = -1524
# Row 58
$ws.Range("H58").Value = 2306.875
$ws.Range("I58").Value = 2001.5714
$ws.Range("K58").Value = 2001.5714
$ws.Range("M58").Value = -1798.5714
# Row 136
$ws.Range("H136").Value = 2306.875
$ws.Range("I136").Value = 2001.5714
$ws.Range("K136").Value = 6004.7142
$ws.Range("M136").Value = -3454.7142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 322.5
$ws.Range("I6").Value = 322.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 967.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -854.5
$ws.Range("N6").ClearContents()
# Row 33
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("K33").Value = 12000
$ws.Range("M33").Value = -11717
# Row 34
$ws.Range("H34").Value = 546.6667
$ws.Range("I34").Value = 550
$ws.Range("J34").Value = 545
$ws.Range("K34").Value = 1650
$ws.Range("L34").Value = 1635
$ws.Range("M34").Value = -1566
$ws.Range("N34").Value = -1803
# Row 40
$ws.Range("H40").Value = 77.125
$ws.Range("J40").Value = 99
$ws.Range("L40").Value = 396
$ws.Range("N40").Value = -534
# Row 69
$ws.Range("H69").Value = 2933.3333
$ws.Range("J69").Value = 2933.3333
$ws.Range("L69").Value = 8799.999899999999
$ws.Range("N69").Value = -10421.9999
# Row 72
$ws.Range("H72").Value = 2933.3333
$ws.Range("J72").Value = 2933.3333
$ws.Range("L72").Value = 26399.9997
$ws.Range("N72").Value = -34511.9997
# Row 98
$ws.Range("H98").Value = 4003
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 4003
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 12009
$ws.Range("N98").Value = -15005
$ws.Range("M98").ClearContents()
# Row 107
$ws.Range("H107").Value = 2003.25
$ws.Range("J107").Value = 2003.3334
$ws.Range("L107").Value = 6010.0002
$ws.Range("N107").Value = -9850.0002
# Row 109
$ws.Range("H109").Value = 999.3333
$ws.Range("I109").Value = 269.2
$ws.Range("K109").Value = 807.5999999999999
$ws.Range("M109").Value = 232.4000000000001
# Row 113
$ws.Range("H113").Value = 1697.3
$ws.Range("J113").Value = 1735.7142
$ws.Range("L113").Value = 5207.142599999999
$ws.Range("N113").Value = -9547.142599999999
# Row 115
$ws.Range("H115").Value = 3514.2856
$ws.Range("J115").Value = 3983.3333
$ws.Range("L115").Value = 11949.9999
$ws.Range("N115").Value = -14299.9999
# Row 131
$ws.Range("H131").Value = 437087.44
$ws.Range("I131").Value = 1137.2
$ws.Range("K131").Value = 3411.6
$ws.Range("M131").Value = 1628.4
# Row 132
$ws.Range("H132").Value = 4111.375
$ws.Range("I132").Value = 3973.75
$ws.Range("J132").Value = 4249
$ws.Range("K132").Value = 35763.75
$ws.Range("L132").Value = 38241
$ws.Range("M132").Value = -33233.75
$ws.Range("N132").Value = -43301

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 510
$ws.Range("I97").Value = 510
$ws.Range("K97").Value = 510
$ws.Range("M97").Value = -14
# Row 132
$ws.Range("H132").Value = 4180.9
$ws.Range("I132").Value = 4145.4443
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 12436.3329
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -9906.332900000001
$ws.Range("N132").Value = -18560

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3475.5
$ws.Range("I136").Value = 1004
$ws.Range("J136").Value = 4299.3335
$ws.Range("K136").Value = 3012
$ws.Range("L136").Value = 12898.0005
$ws.Range("M136").Value = -462
$ws.Range("N136").Value = -17998.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 478875
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 478875
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 478875
$ws.Range("N2").Value = -479099
$ws.Range("M2").ClearContents()
# Row 3
$ws.Range("H3").Value = 46000
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 90000
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 90000
$ws.Range("M3").Value = -1886
$ws.Range("N3").Value = -90228
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 52
$ws.Range("H52").Value = 29832.666
$ws.Range("J52").Value = 29832.666
$ws.Range("L52").Value = 29832.666
$ws.Range("N52").Value = -30284.666
# Row 70
$ws.Range("H70").Value = 84800
$ws.Range("J70").Value = 84800
$ws.Range("L70").Value = 84800
$ws.Range("N70").Value = -85430
# Row 73
$ws.Range("H73").Value = 84800
$ws.Range("J73").Value = 84800
$ws.Range("L73").Value = 84800
$ws.Range("N73").Value = -86984
# Row 75
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 74801.25
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 74801.25
$ws.Range("N75").Value = -76673.25
$ws.Range("M75").ClearContents()
# Row 78
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 74801.25
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 224403.75
$ws.Range("N78").Value = -233763.75
$ws.Range("M78").ClearContents()
# Row 126
$ws.Range("H126").Value = 4661.657
$ws.Range("I126").Value = 4622.4546
$ws.Range("J126").Value = 4728
$ws.Range("K126").Value = 13867.3638
$ws.Range("L126").Value = 14184
$ws.Range("M126").Value = -11397.3638
$ws.Range("N126").Value = -19124
